$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.729.91'
$ws.Range("E2").Value = '  +4.25%  '
$ws.Range("D3").Value = '2.262.29'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'300.78"
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = "'100.32"
$ws.Range("E6").Value = '  +6.67%  '
$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = "'0.511"
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = "'35.57"
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("D11").Value = "'0.0780"
$ws.Range("E11").Value = '  -1.14%  '
$ws.Range("D12").Value = "'7.16"
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '2.605.00'
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").Value = '2.256.60'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").Value = "'13.58"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '46.686.78'
$ws.Range("E17").Value = '  +4.31%  '
$ws.Range("D18").Value = "'0.793"
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").Value = "'12.70"
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").Value = '0.0₃0925'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = "'5.86"
$ws.Range("E21").Value = '  -3.07%  '
$ws.Range("D22").Value = "'65.15"
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").Value = "'248.46"
$ws.Range("E23").Value = '  +4.38%  '
$ws.Range("D24").Value = "'2.81"
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").Value = "'1.87"
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = "'42.41"
$ws.Range("E27").Value = '  +2.69%  '
$ws.Range("D28").Value = "'2.24"
$ws.Range("E28").Value = '  -2.56%  '
$ws.Range("D29").Value = "'9.67"
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = "'19.84"
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").Value = "'2.75"
$ws.Range("E31").Value = '  +7.95%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'5.42"
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").Value = "'145.02"
$ws.Range("E33").Value = '  -4.66%  '
$ws.Range("D34").Value = "'0.0771"
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("D35").Value = "'3.19"
$ws.Range("E35").Value = '  +8.21%  '
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = '  +9.40%  '
$ws.Range("E37").Value = '  -1.67%  '
$ws.Range("D38").Value = "'16.02"
$ws.Range("E38").Value = '  +18.25%  '
$ws.Range("D39").Value = "'1.71"
$ws.Range("E39").Value = '  -2.91%  '
$ws.Range("D40").Value = "'3.86"
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("D41").Value = "'0.0297"
$ws.Range("E41").Value = '  -4.08%  '
$ws.Range("D42").Value = "'3.20"
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").Value = "'1.97"
$ws.Range("E44").Value = '  +3.56%  '
$ws.Range("D45").Value = '1.794.38'
$ws.Range("E45").Value = '  +3.60%  '
$ws.Range("D46").Value = "'91.31"
$ws.Range("E46").Value = '  +20.21%  '
$ws.Range("D47").Value = "'0.188"
$ws.Range("E47").Value = '  -2.96%  '
$ws.Range("D48").Value = "'71.49"
$ws.Range("E48").Value = '  +3.24%  '
$ws.Range("D49").Value = "'4.80"
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("D50").Value = "'93.68"
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("D51").Value = '2.484.29'
$ws.Range("E51").Value = '  -0.08%  '
